$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (workbook.xml <sheet name="...">)
$ws.Name = "BOM_Board1_PCB1_2024-11-14"

# Only the cells that actually change value are touched below, so untouched
# cells (including the always-blank "Value" column F) are left exactly as
# they were.

# Row 2 becomes the new "Test-Point" line (No. 1), designator gains "5V,"
$ws.Range("C2").Value = "Test-Point"
$ws.Range("D2").Value = "5V,DO,RGB"
$ws.Range("E2").Value = "Test-Point-0.5mm"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# Row 3 becomes the "1uF" line (No. 2)
$ws.Range("C3").Value = "1uF"
$ws.Range("D3").Value = "C1,C2,C3"
$ws.Range("G3").Value = "CS0805KKX7R8BB105"
$ws.Range("H3").Value = "YAGEO(国巨)"
$ws.Range("I3").Value = "C272869"

# Row 4 becomes the "100nF" line (No. 3)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "100nF"
$ws.Range("D4").Value = "C4,C5,C6"
$ws.Range("E4").Value = "C0805"
$ws.Range("G4").Value = "0805B104K101AT"
$ws.Range("H4").Value = "FH(风华)"
$ws.Range("I4").Value = "C3037660"

# Row 5 becomes the "1N4148WS-V-GS18" line (No. 4)
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "1N4148WS-V-GS18"
$ws.Range("D5").Value = "D1"
$ws.Range("E5").Value = "SOD-323_L1.8-W1.3-LS2.7-RD"
$ws.Range("G5").Value = "1N4148WS-V-GS18"
$ws.Range("H5").Value = "VISHAY(威世)"
$ws.Range("I5").Value = "C19078477"

# Row 6 becomes the "SM06B-SRSS-TB(LF)(SN)" line (No. 5)
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "SM06B-SRSS-TB(LF)(SN)"
$ws.Range("D6").Value = "DISP,i2c_EXP,RGB_EXP,U3"
$ws.Range("E6").Value = "CONN-SMD_SM06B-SRSS-TB-LF-SN"
$ws.Range("G6").Value = "SM06B-SRSS-TB(LF)(SN)"
$ws.Range("H6").Value = "JST"
$ws.Range("I6").Value = "C160405"
$ws.Range("J6").Value = "LCSC"

# Row 7 (No. 6, L1) swaps the inductor part for a 220-ohm resistor entry
$ws.Range("C7").Value = "220OHM-0805L"
$ws.Range("G7").Value = "220ohm-0805L"
$ws.Range("H7").Value = "null"
$ws.Range("I7").Value = "C9900020208"

# Row 17 (No. 16) swaps the LDO regulator part
$ws.Range("C17").Value = "MIC5504-3.3YM5-TR"
$ws.Range("D17").Value = "U4"
$ws.Range("E17").Value = "SOT-23-5_L3.0-W1.7-P0.95-LS2.8-TL"
$ws.Range("G17").Value = "MIC5504-3.3YM5-TR"
$ws.Range("H17").Value = "MICROCHIP(美国微芯)"
$ws.Range("I17").Value = "C88419"
